$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B, rows 4 through 73 (naive forecaster bugfix recomputation)
$newValues = @{
    4 = 0.9623985945846414
    5 = -0.2625025240627394
    6 = -2.717991563576632
    7 = -0.8018087606025261
    8 = -2.540735588184205
    9 = 0.4798059474883161
    10 = 2.996743291460049
    11 = 1.766782973262963
    12 = 1.792707332192279
    13 = 1.599999999999994
    14 = 1.983944815439088
    15 = 1.027205282249909
    16 = 1.892417316869597
    17 = 1.727088565964991
    18 = 0.4000000000000057
    19 = 0.5999999999999943
    20 = 1.295244683175738
    21 = 0.5000000000000142
    22 = -0.1007444012410019
    23 = 0.399975276622385
    24 = 2
    25 = 1.1
    26 = 1.292635181922734
    27 = 1.682020243440505
    28 = 1.012497979540356
    29 = 1.634557118349079
    30 = 0.8971548841028039
    31 = 1.102490924221428
    32 = 0.9984536597660991
    33 = 1.495216092286043
    34 = 1.345302727311861
    35 = 1.428495556385869
    36 = 0.7394362573787987
    37 = 1.299999999999983
    38 = 1.233394657233262
    39 = 1.297884859680252
    40 = 1.637264503497377
    41 = 1.331979115623398
    42 = 1.795921598270084
    43 = 1.327129713066284
    44 = 1.15883627336575
    45 = 1.076058203620576
    46 = 0.3
    47 = 2.40590956953757
    48 = -0.4
    49 = 0.8963263664365542
    50 = 1.055400706275506
    51 = -0.4303992348575321
    52 = -14.5
    53 = 7.96955251685678
    54 = 2.117788110998191
    55 = -0.9861240056009706
    56 = 2.876944405321424
    57 = 0.4364757668776207
    58 = 0.6542354095451515
    59 = -1.358640149334988
    60 = 1.767346889326234
    61 = 0.1973819540654631
    62 = -2.016521230865749
    63 = -1.429923541452922
    64 = 0.1801827438520291
    65 = 0.6692453970872521
    66 = -1.33243152085096
    67 = -0.09168777270478756
    68 = 0.8481739611978583
    69 = -0.2604183589432552
    70 = 0.3093370292089048
    71 = -0.03180050048325711
    72 = 0.04735640278761366
    73 = 0.2461857363876589
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}

# Rows 74:82 are no longer part of the series - remove them entirely
$ws.Rows("74:82").Delete()
